$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-10-02 01:13:34"
}
